$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New cell values, entered in the order that reproduces the shared-strings table order
$ws.Range("C15").Value = "Under Frame"
$ws.Range("C17").Value = "Bracing à découper"
$ws.Range("C16").Value = "Frame Top View"
$ws.Range("B17").Value = "20_G_FR_A0100_#07"
$ws.Range("B18").Value = "20_G_FR_A0100_201"
$ws.Range("B19").Value = "20_G_FR_A0100_202"
$ws.Range("C18").Value = "Front Hoop"
$ws.Range("C19").Value = "Main Hoop"
$ws.Range("E18").Value = "MJT"
$ws.Range("E19").Value = "MJT"
$ws.Range("C14").Value = "Cellule Arrière"

# Column D updates ("Oui")
$ws.Range("D12").Value = "Oui"
$ws.Range("D14").Value = "Oui"
$ws.Range("D16").Value = "Oui"
$ws.Range("D17").Value = "Oui"
$ws.Range("D18").Value = "Oui"
$ws.Range("D19").Value = "Oui"

# Remove the border style on C18 (Excel reports s="1" there instead of s="2")
$ws.Range("C18").Borders.LineStyle = -4142

# Window view changes
$ws.Range("D15").Select()
$excel.ActiveWindow.ScrollRow = 10

$wb.Save()
